# Test Data changes - Companies module - 4th Dec 2023
# Update the Users sheet: replace the user name "Drew Koecher" with "Ayati Arvind"
# and leave the selection on cell H6 as the active cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

$ws.Range("A2").Value = "Ayati Arvind"

$ws.Range("H6").Select()
